$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'41.798.00"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.02%  "

# Row 3
$ws.Range("D3").Value = "'2.217.45"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.29%  "

# Row 4
$ws.Range("E4").Value = "  -0.06%  "

# Row 5
$ws.Range("D5").Value = "'252.36"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.21%  "

# Row 6
$ws.Range("D6").Value = "'0.630"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.33%  "

# Row 7
$ws.Range("D7").Value = "'71.19"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.60%  "

# Row 8
$ws.Range("E8").Value = "  +0.06%  "

# Row 9
$ws.Range("D9").Value = "'0.595"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +7.20%  "

# Row 10
$ws.Range("D10").Value = "'40.66"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +11.89%  "

# Row 11
$ws.Range("D11").Value = "'0.0965"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.48%  "

# Row 12
$ws.Range("D12").Value = "'58.22"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.18%  "

# Row 13
$ws.Range("D13").Value = "'7.22"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +5.92%  "

# Row 14
$ws.Range("E14").Value = "  -1.45%  "

# Row 15
$ws.Range("D15").Value = "'2.546.87"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.40%  "

# Row 16
$ws.Range("D16").Value = "'14.97"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.41%  "

# Row 17
$ws.Range("D17").Value = "'0.876"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.97%  "

# Row 18
$ws.Range("D18").Value = "'2.217.98"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.62%  "

# Row 19
$ws.Range("D19").Value = "'41.796.73"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.01%  "

# Row 20
$ws.Range("D20").Value = "'0.0₃0961"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.13%  "

# Row 21
$ws.Range("D21").Value = "'6.21"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.87%  "

# Row 22
$ws.Range("D22").Value = "'72.61"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.09%  "

# Row 23
$ws.Range("D23").Value = "'235.25"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.47%  "

# Row 24
$ws.Range("D24").Value = "'2.06"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.94%  "

# Row 25
$ws.Range("D25").Value = "'4.10"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +11.75%  "

# Row 26
$ws.Range("E26").Value = "  +0.12%  "

# Row 27
$ws.Range("D27").Value = "'2.54"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.35%  "

# Row 28
$ws.Range("D28").Value = "'11.13"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +10.82%  "

# Row 29
$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").Value = "'2.21"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.39%  "

# Row 30
$ws.Range("B30").Value = "Monero"
$ws.Range("C30").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D30").Value = "'170.01"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.73%  "

# Row 31
$ws.Range("D31").Value = "'20.76"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.92%  "

# Row 32
$ws.Range("E32").Value = "  -0.81%  "

# Row 33
$ws.Range("B33").Value = "Hedera"
$ws.Range("C33").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D33").Value = "'0.0740"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.81%  "

# Row 34
$ws.Range("B34").Value = "InternetComputer(DFINITY)"
$ws.Range("C34").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D34").Value = "'5.51"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.67%  "

# Row 35
$ws.Range("B35").Value = "Stellar"
$ws.Range("C35").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D35").Value = "'0.122"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.52%  "

# Row 36
$ws.Range("D36").Value = "'26.60"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +16.68%  "

# Row 37
$ws.Range("E37").Value = "  -0.62%  "

# Row 38
$ws.Range("D38").Value = "'4.02"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +4.07%  "

# Row 39
$ws.Range("D39").Value = "'0.0306"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +8.45%  "

# Row 40
$ws.Range("D40").Value = "'2.28"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.85%  "

# Row 41
$ws.Range("D41").Value = "'12.75"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +26.52%  "

# Row 42
$ws.Range("D42").Value = "'5.93"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.02%  "

# Row 43
$ws.Range("D43").Value = "'64.18"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.49%  "

# Row 44
$ws.Range("E44").Value = "  +6.35%  "

# Row 45
$ws.Range("B45").Value = "FraxShare"
$ws.Range("C45").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D45").Value = "'8.73"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -7.13%  "

# Row 46
$ws.Range("E46").Value = "  -0.90%  "

# Row 47
$ws.Range("B47").Value = "FTXToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D47").Value = "'4.71"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -5.09%  "

# Row 48
$ws.Range("B48").Value = "SynthetixNetwork"
$ws.Range("C48").Value = "https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx"
$ws.Range("D48").Value = "'4.68"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.40%  "

# Row 49
$ws.Range("E49").Value = "  +0.00%  "

# Row 50
$ws.Range("D50").Value = "'1.17"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +4.60%  "

# Row 51
$ws.Range("D51").Value = "'2.37"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.63%  "
